$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets("ALC")
$ws.Range("H74").Value = 3000
$ws.Range("I74").Value = 3000
$ws.Range("K74").Value = 3000
$ws.Range("M74").Value = -2064
$ws.Range("H76").Value = 2000
$ws.Range("I76").Value = 2000
$ws.Range("K76").Value = 2000
$ws.Range("M76").Value = -1685
$ws.Range("H77").Value = 3000
$ws.Range("I77").Value = 3000
$ws.Range("K77").Value = 15000
$ws.Range("M77").Value = -10320
$ws.Range("H79").Value = 2000
$ws.Range("I79").Value = 2000
$ws.Range("K79").Value = 2000
$ws.Range("M79").Value = -908
$ws.Range("H92").Value = 433.85715
$ws.Range("I92").Value = 207.6
$ws.Range("K92").Value = 207.6
$ws.Range("M92").Value = 1040.4
$ws.Range("H98").Value = 8980.75
$ws.Range("I98").Value = 8980.75
$ws.Range("K98").Value = 8980.75
$ws.Range("M98").Value = -7482.75
$ws.Range("H107").Value = 1289.8462
$ws.Range("I107").Value = 818.7
$ws.Range("K107").Value = 818.7
$ws.Range("M107").Value = 1101.3
$ws.Range("H122").Value = 8980.75
$ws.Range("I122").Value = 8980.75
$ws.Range("K122").Value = 26942.25
$ws.Range("M122").Value = -24492.25
$ws.Range("H125").Value = 27500
$ws.Range("I125").Value = 25000
$ws.Range("J125").Value = 30000
$ws.Range("K125").Value = 225000
$ws.Range("L125").Value = 270000
$ws.Range("M125").Value = -222540
$ws.Range("N125").Value = -274920
$ws.Range("H135").Value = 1797.5555
$ws.Range("I135").Value = 1141.6666
$ws.Range("K135").Value = 10274.9994
$ws.Range("M135").Value = -7739.999400000001
$ws.Range("H137").Value = 1364
$ws.Range("I137").Value = 1415.1666
$ws.Range("K137").Value = 4245.4998
$ws.Range("M137").Value = -1695.4998

# --- ARM ---
$ws = $wb.Worksheets("ARM")
$ws.Range("H61").Value = 2000
$ws.Range("I61").Value = 2000
$ws.Range("K61").Value = 2000
$ws.Range("M61").Value = -1788
$ws.Range("H63").Value = 2666.6667
$ws.Range("J63").Value = 1500
$ws.Range("L63").Value = 1500
$ws.Range("N63").Value = -2872
$ws.Range("H66").Value = 2666.6667
$ws.Range("J66").Value = 1500
$ws.Range("L66").Value = 7500
$ws.Range("N66").Value = -14364
$ws.Range("H74").Value = 3321
$ws.Range("I74").Value = 3574.8333
$ws.Range("K74").Value = 3574.8333
$ws.Range("M74").Value = -2700.8333
$ws.Range("H77").Value = 3321
$ws.Range("I77").Value = 3574.8333
$ws.Range("K77").Value = 17874.1665
$ws.Range("M77").Value = -13506.1665
$ws.Range("H97").Value = 1030.5834
$ws.Range("I97").Value = 336.7
$ws.Range("K97").Value = 336.7
$ws.Range("M97").Value = 159.3
$ws.Range("H124").Value = 43464
$ws.Range("J124").Value = 43464
$ws.Range("L124").Value = 43464
$ws.Range("N124").Value = -53284
$ws.Range("H136").Value = 2000
$ws.Range("I136").Value = 2000
$ws.Range("K136").Value = 6000
$ws.Range("M136").Value = -3450

# --- CRP ---
$ws = $wb.Worksheets("CRP")
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 250
$ws.Range("K4").Value = 250
$ws.Range("M4").Value = -138
$ws.Range("H31").Value = 7360
$ws.Range("J31").Value = 2647.5
$ws.Range("L31").Value = 2647.5
$ws.Range("N31").Value = -3237.5
$ws.Range("H34").Value = 7360
$ws.Range("J34").Value = 2647.5
$ws.Range("L34").Value = 2647.5
$ws.Range("N34").Value = -3051.5
$ws.Range("H132").Value = 3194.4546
$ws.Range("I132").Value = 3208.25
$ws.Range("J132").Value = 3157.6667
$ws.Range("K132").Value = 9624.75
$ws.Range("L132").Value = 9473.000100000001
$ws.Range("M132").Value = -7094.75
$ws.Range("N132").Value = -14533.0001
$ws.Range("H134").Value = 2049.5
$ws.Range("I134").Value = 2110.5557
$ws.Range("J134").Value = 1500
$ws.Range("K134").Value = 6331.6671
$ws.Range("L134").Value = 4500
$ws.Range("M134").Value = -3796.6671
$ws.Range("N134").Value = -9570

# --- CUL ---
$ws = $wb.Worksheets("CUL")
$ws.Range("H68").Value = 1028
$ws.Range("J68").Value = 1003
$ws.Range("L68").Value = 3009
$ws.Range("N68").Value = -4631
$ws.Range("H71").Value = 1028
$ws.Range("J71").Value = 1003
$ws.Range("L71").Value = 9027
$ws.Range("N71").Value = -17139

# --- GSM ---
$ws = $wb.Worksheets("GSM")
$ws.Range("H21").Value = 564099.4399999999
$ws.Range("I21").Value = 1005679
$ws.Range("J21").Value = 12125
$ws.Range("K21").Value = 1005679
$ws.Range("L21").Value = 12125
$ws.Range("M21").Value = -1005506
$ws.Range("N21").Value = -12471
$ws.Range("H30").Value = 564099.4399999999
$ws.Range("I30").Value = 1005679
$ws.Range("J30").Value = 12125
$ws.Range("K30").Value = 1005679
$ws.Range("L30").Value = 12125
$ws.Range("M30").Value = -1005574
$ws.Range("N30").Value = -12335
$ws.Range("H70").Value = 9331.666999999999
$ws.Range("I70").Value = 9331.666999999999
$ws.Range("K70").Value = 9331.666999999999
$ws.Range("M70").Value = -9061.666999999999
$ws.Range("H73").Value = 9331.666999999999
$ws.Range("I73").Value = 9331.666999999999
$ws.Range("K73").Value = 9331.666999999999
$ws.Range("M73").Value = -8395.666999999999
$ws.Range("H92").Value = 20000
$ws.Range("J92").Value = 20000
$ws.Range("L92").Value = 20000
$ws.Range("N92").Value = -23744
$ws.Range("H99").Value = 11249.75
$ws.Range("I99").Value = 4999.5
$ws.Range("K99").Value = 4999.5
$ws.Range("M99").Value = -2753.5
$ws.Range("H113").Value = 666.3333
$ws.Range("I113").Value = 600
$ws.Range("K113").Value = 600
$ws.Range("M113").Value = 1570
$ws.Range("H122").Value = 13543.125
$ws.Range("I122").Value = 1899
$ws.Range("K122").Value = 5697
$ws.Range("M122").Value = -3247

# --- LTW ---
$ws = $wb.Worksheets("LTW")
$ws.Range("H22").Value = 600
$ws.Range("I22").Value = 600
$ws.Range("K22").Value = 600
$ws.Range("M22").Value = -305
$ws.Range("H27").Value = 600
$ws.Range("I27").Value = 600
$ws.Range("K27").Value = 600
$ws.Range("M27").Value = -493
$ws.Range("H40").Value = 6050.8
$ws.Range("J40").Value = 7502.5
$ws.Range("L40").Value = 7502.5
$ws.Range("N40").Value = -7774.5
$ws.Range("H93").Value = 1475.5
$ws.Range("J93").Value = 838
$ws.Range("L93").Value = 838
$ws.Range("N93").Value = -3334
$ws.Range("H122").Value = 7649.8335
$ws.Range("I122").Value = 7649.8335
$ws.Range("K122").Value = 22949.5005
$ws.Range("M122").Value = -20499.5005

# --- WVR ---
$ws = $wb.Worksheets("WVR")
$ws.Range("H74").Value = 49999
$ws.Range("I74").Value = 49998.5
$ws.Range("K74").Value = 49998.5
$ws.Range("M74").Value = -49062.5
$ws.Range("H77").Value = 49999
$ws.Range("I77").Value = 49998.5
$ws.Range("K77").Value = 149995.5
$ws.Range("M77").Value = -145315.5
$ws.Range("H122").Value = 224925.11
$ws.Range("I122").Value = 287260.84
$ws.Range("J122").Value = 6750
$ws.Range("K122").Value = 861782.52
$ws.Range("L122").Value = 20250
$ws.Range("M122").Value = -859332.52
$ws.Range("N122").Value = -25150
